# Update cryptos list (coinranking.com price snapshot refresh).
# Mirrors the per-row Price (D) / Volume(1h) (E) updates from the source diff,
# plus a full swap of the Aptos / TheSandbox rows (40 <-> 41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "29.745.02"
$ws.Range("E2").Value = "  -2.56%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.098.54"
$ws.Range("E3").Value = "  -1.93%  "

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.21"
$ws.Range("E5").Value = "  -2.29%  "

# Row 6: USDC
$ws.Range("E6").Value = "  +0.17%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5171"
$ws.Range("E7").Value = "  -1.66%  "

# Row 8: Cardano
$ws.Range("E8").Value = "  -3.82%  "

# Row 9: OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.06"
$ws.Range("E9").Value = "  -1.10%  "

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09304"
$ws.Range("E10").Value = "  +1.29%  "

# Row 11: Polygon
$ws.Range("E11").Value = "  -2.41%  "

# Row 12: Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.97"
$ws.Range("E12").Value = "  -2.38%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "2.107.31"
$ws.Range("E13").Value = "  -0.79%  "

# Row 14: Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.268"
$ws.Range("E14").Value = "  +1.12%  "

# Row 15: Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.777"
$ws.Range("E15").Value = "  -1.82%  "

# Row 16: Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.68"
$ws.Range("E16").Value = "  -1.89%  "

# Row 17: ShibaInu
$ws.Range("E17").Value = "  -1.32%  "

# Row 19: Avalanche
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.87"
$ws.Range("E19").Value = "  +1.36%  "

# Row 20: TRON
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06648"
$ws.Range("E20").Value = "  -1.15%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.208"
$ws.Range("E22").Value = "  -2.73%  "

# Row 23: WrappedBTC
$ws.Range("D23").Value = "29.775.25"
$ws.Range("E23").Value = "  -2.73%  "

# Row 24: Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.52"
$ws.Range("E24").Value = "  -2.67%  "

# Row 25: Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.315"
$ws.Range("E25").Value = "  -2.60%  "

# Row 26: WrappedliquidstakedEther2.0
$ws.Range("D26").Value = "2.348.07"
$ws.Range("E26").Value = "  -1.34%  "

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.97"
$ws.Range("E27").Value = "  -2.47%  "

# Row 28: LidoDAOToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.523"
$ws.Range("E28").Value = "  -3.19%  "

# Row 29: Monero
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.43"
$ws.Range("E29").Value = "  -2.27%  "

# Row 30: BitcoinCash
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.15"
$ws.Range("E30").Value = "  -2.08%  "

# Row 31: ImmutableX
$ws.Range("E31").Value = "  -6.68%  "

# Row 32: Stellar
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1050"
$ws.Range("E32").Value = "  -3.06%  "

# Row 33: ARBITRUM
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.654"
$ws.Range("E33").Value = "  -3.91%  "

# Row 34: Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.186"
$ws.Range("E34").Value = "  -3.58%  "

# Row 35: HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.943"
$ws.Range("E35").Value = "  -2.40%  "

# Row 36: InternetComputer(DFINITY)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.284"
$ws.Range("E36").Value = "  +1.80%  "

# Row 37: FraxShare
$ws.Range("E37").Value = "  -1.92%  "

# Row 38: VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02576"
$ws.Range("E38").Value = "  -2.72%  "

# Row 39: Hedera
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06732"
$ws.Range("E39").Value = "  -3.78%  "

# Row 40: TheSandbox
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6928"
$ws.Range("E40").Value = "  -1.12%  "

# Row 41: Aptos
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.47"
$ws.Range("E41").Value = "  -2.18%  "

# Row 42: Algorand
$ws.Range("E42").Value = "  -5.31%  "

# Row 43: TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.316"
$ws.Range("E43").Value = "  +3.09%  "

# Row 44: Decentraland
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6820"
$ws.Range("E44").Value = "  +4.38%  "

# Row 45: EnergySwap
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.40"
$ws.Range("E45").Value = "  -3.14%  "

# Row 46: NEARProtocol
$ws.Range("E46").Value = "  -1.34%  "

# Row 47: PancakeSwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.626"
$ws.Range("E47").Value = "  -2.61%  "

# Row 48: BabyDogeCoin
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000356"
$ws.Range("E48").Value = "  -5.56%  "

# Row 49: EOS
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.221"
$ws.Range("E49").Value = "  -2.56%  "

# Row 50: Aave
$ws.Range("E50").Value = "  -2.06%  "

# Row 51: WEMIXTOKEN
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.161"
$ws.Range("E51").Value = "  -2.10%  "
